$wb = $excel.ActiveWorkbook

# Each entry: Sheet name, cell address, new value ($null clears the cell)
$edits = @(
    @("ALC", "H4", 134.875),
    @("ALC", "I4", 125.57143),
    @("ALC", "K4", 125.57143),
    @("ALC", "M4", -11.57143000000001),
    @("ALC", "H62", 3638),
    @("ALC", "I62", 3451.2222),
    @("ALC", "K62", 3451.2222),
    @("ALC", "M62", -2827.2222),
    @("ALC", "H65", 3638),
    @("ALC", "I65", 3451.2222),
    @("ALC", "K65", 17256.111),
    @("ALC", "M65", -14136.111),
    @("ALC", "H141", 2446.7646),
    @("ALC", "I141", 2446.7646),
    @("ALC", "K141", 7340.293799999999),
    @("ALC", "M141", -2160.293799999999),
    @("ARM", "H34", 30000),
    @("ARM", "J34", 30000),
    @("ARM", "L34", 30000),
    @("ARM", "N34", -30542),
    @("ARM", "H37", 16999),
    @("ARM", "J37", 19998.75),
    @("ARM", "L37", 19998.75),
    @("ARM", "N37", -20544.75),
    @("ARM", "H97", 534.5714),
    @("ARM", "I97", 552.61536),
    @("ARM", "J97", 300),
    @("ARM", "K97", 552.61536),
    @("ARM", "L97", 300),
    @("ARM", "M97", -56.61536000000001),
    @("ARM", "N97", -1292),
    @("ARM", "H102", 1482.7),
    @("ARM", "I102", 1004),
    @("ARM", "K102", 1004),
    @("ARM", "M102", 618),
    @("BSM", "H86", 2156.3794),
    @("BSM", "I86", 1452.6111),
    @("BSM", "J86", 3308),
    @("BSM", "K86", 1452.6111),
    @("BSM", "L86", 3308),
    @("BSM", "M86", -329.6111000000001),
    @("BSM", "N86", -5554),
    @("BSM", "H89", 2156.3794),
    @("BSM", "I89", 1452.6111),
    @("BSM", "J89", 3308),
    @("BSM", "K89", 7263.0555),
    @("BSM", "L89", 16540),
    @("BSM", "M89", -1647.0555),
    @("BSM", "N89", -27772),
    @("BSM", "H99", 2495.2144),
    @("BSM", "I99", 2403),
    @("BSM", "K99", 2403),
    @("BSM", "M99", -905),
    @("CRP", "H16", 2013),
    @("CRP", "J16", 2013),
    @("CRP", "L16", 2013),
    @("CRP", "N16", -2587),
    @("CRP", "H31", 4131),
    @("CRP", "I31", 2440.6667),
    @("CRP", "K31", 2440.6667),
    @("CRP", "M31", -2145.6667),
    @("CRP", "H34", 4131),
    @("CRP", "I34", 2440.6667),
    @("CRP", "K34", 2440.6667),
    @("CRP", "M34", -2238.6667),
    @("CRP", "H86", 9160.223),
    @("CRP", "I86", 9129.200000000001),
    @("CRP", "K86", 9129.200000000001),
    @("CRP", "M86", -8006.200000000001),
    @("CRP", "H89", 9160.223),
    @("CRP", "I89", 9129.200000000001),
    @("CRP", "K89", 45646),
    @("CRP", "M89", -40030),
    @("CRP", "H113", 2013),
    @("CRP", "J113", 2013),
    @("CRP", "L113", 2013),
    @("CRP", "N113", -6353),
    @("CRP", "H132", 2179.2),
    @("CRP", "I132", 2437.5386),
    @("CRP", "J132", 500),
    @("CRP", "K132", 7312.6158),
    @("CRP", "L132", 1500),
    @("CRP", "M132", -4782.6158),
    @("CRP", "N132", -6560),
    @("CUL", "H113", 910.8),
    @("CUL", "J113", 1034.6154),
    @("CUL", "L113", 3103.8462),
    @("CUL", "N113", -7443.8462),
    @("CUL", "H122", 2121),
    @("CUL", "J122", 2121),
    @("CUL", "L122", 19089),
    @("CUL", "N122", -23989),
    @("GSM", "H11", 685000),
    @("GSM", "I11", 850000),
    @("GSM", "J11", 25000),
    @("GSM", "K11", 850000),
    @("GSM", "L11", 25000),
    @("GSM", "M11", -849861),
    @("GSM", "N11", -25278),
    @("GSM", "H122", 3002.48),
    @("GSM", "I122", 2875.7646),
    @("GSM", "J122", 3271.75),
    @("GSM", "K122", 8627.293799999999),
    @("GSM", "L122", 9815.25),
    @("GSM", "M122", -6177.293799999999),
    @("GSM", "N122", -14715.25),
    @("GSM", "H132", 2998.6365),
    @("GSM", "I132", 2998.6365),
    @("GSM", "K132", 8995.9095),
    @("GSM", "M132", -6465.9095),
    @("LTW", "H61", 8701.6),
    @("LTW", "I61", 7833.3335),
    @("LTW", "J61", 10004),
    @("LTW", "K61", 7833.3335),
    @("LTW", "L61", 10004),
    @("LTW", "M61", -7631.3335),
    @("LTW", "N61", -10408),
    @("LTW", "H100", 2055.4285),
    @("LTW", "I100", 877.6),
    @("LTW", "K100", 877.6),
    @("LTW", "M100", -336.6),
    @("LTW", "H113", 8701.6),
    @("LTW", "I113", 7833.3335),
    @("LTW", "J113", 10004),
    @("LTW", "K113", 7833.3335),
    @("LTW", "L113", 10004),
    @("LTW", "M113", -5663.3335),
    @("LTW", "N113", -14344),
    @("WVR", "H18", 499),
    @("WVR", "I18", 0),
    @("WVR", "K18", 0),
    @("WVR", "M18", $null),
    @("WVR", "H81", 2098.6),
    @("WVR", "I81", 2098.6),
    @("WVR", "J81", 0),
    @("WVR", "K81", 4197.2),
    @("WVR", "L81", 0),
    @("WVR", "M81", -3136.2),
    @("WVR", "N81", $null),
    @("WVR", "H84", 2098.6),
    @("WVR", "I84", 2098.6),
    @("WVR", "J84", 0),
    @("WVR", "K84", 23061.428),
    @("WVR", "L84", 0),
    @("WVR", "M84", -15682),
    @("WVR", "N84", $null),
    @("WVR", "H113", 418.83334),
    @("WVR", "I113", 418.83334),
    @("WVR", "J113", 0),
    @("WVR", "K113", 1256.50002),
    @("WVR", "L113", 0),
    @("WVR", "M113", 913.4999800000001),
    @("WVR", "N113", $null),
    @("WVR", "H122", 3797.8572),
    @("WVR", "J122", 4999.5),
    @("WVR", "L122", 14998.5),
    @("WVR", "N122", -19898.5),
    @("WVR", "H136", 3242.1667),
    @("WVR", "I136", 3188.0356),
    @("WVR", "K136", 9564.106800000001),
    @("WVR", "M136", -7014.106800000001)
)

foreach ($e in $edits) {
    $ws = $wb.Worksheets.Item($e[0])
    $ws.Range($e[1]).Value = $e[2]
}
